$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1895.8572
$ws.Range("I98").Value = 1334.7391
$ws.Range("J98").Value = 4477
$ws.Range("K98").Value = 1334.7391
$ws.Range("L98").Value = 4477
$ws.Range("M98").Value = 163.2609
$ws.Range("N98").Value = -7473

$ws.Range("H122").Value = 1895.8572
$ws.Range("I122").Value = 1334.7391
$ws.Range("J122").Value = 4477
$ws.Range("K122").Value = 4004.2173
$ws.Range("L122").Value = 13431
$ws.Range("M122").Value = -1554.2173
$ws.Range("N122").Value = -18331

$ws.Range("H137").Value = 1334.9286
$ws.Range("I137").Value = 1307.1666
$ws.Range("K137").Value = 3921.4998
$ws.Range("M137").Value = -1371.4998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1968.4166
$ws.Range("I61").Value = 1713.375
$ws.Range("J61").Value = 2478.5
$ws.Range("K61").Value = 1713.375
$ws.Range("L61").Value = 2478.5
$ws.Range("M61").Value = -1501.375
$ws.Range("N61").Value = -2902.5

$ws.Range("H74").Value = 2740
$ws.Range("I74").Value = 970.6667
$ws.Range("K74").Value = 970.6667
$ws.Range("M74").Value = -96.66669999999999

$ws.Range("H77").Value = 2740
$ws.Range("I77").Value = 970.6667
$ws.Range("K77").Value = 4853.3335
$ws.Range("M77").Value = -485.3334999999997

$ws.Range("H88").Value = 2943.4285
$ws.Range("J88").Value = 2943.4285
$ws.Range("L88").Value = 2943.4285
$ws.Range("N88").Value = -3755.4285

$ws.Range("H91").Value = 2943.4285
$ws.Range("J91").Value = 2943.4285
$ws.Range("L91").Value = 2943.4285
$ws.Range("N91").Value = -5751.4285

$ws.Range("H110").Value = 2067.8
$ws.Range("I110").Value = 1438.8334
$ws.Range("K110").Value = 1438.8334
$ws.Range("M110").Value = 606.1666

$ws.Range("H111").Value = 40644
$ws.Range("J111").Value = 40644
$ws.Range("L111").Value = 40644
$ws.Range("N111").Value = -48824

$ws.Range("H132").Value = 2783.647
$ws.Range("I132").Value = 2768.6667
$ws.Range("J132").Value = 2819.6
$ws.Range("K132").Value = 8306.000100000001
$ws.Range("L132").Value = 8458.799999999999
$ws.Range("M132").Value = -5776.000100000001
$ws.Range("N132").Value = -13518.8

$ws.Range("H136").Value = 1968.4166
$ws.Range("I136").Value = 1713.375
$ws.Range("J136").Value = 2478.5
$ws.Range("K136").Value = 5140.125
$ws.Range("L136").Value = 7435.5
$ws.Range("M136").Value = -2590.125
$ws.Range("N136").Value = -12535.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2502.7778
$ws.Range("I31").Value = 1482.5714
$ws.Range("J31").Value = 2859.85
$ws.Range("K31").Value = 1482.5714
$ws.Range("L31").Value = 2859.85
$ws.Range("M31").Value = -1187.5714
$ws.Range("N31").Value = -3449.85

$ws.Range("H34").Value = 2502.7778
$ws.Range("I34").Value = 1482.5714
$ws.Range("J34").Value = 2859.85
$ws.Range("K34").Value = 1482.5714
$ws.Range("L34").Value = 2859.85
$ws.Range("M34").Value = -1280.5714
$ws.Range("N34").Value = -3263.85

$ws.Range("H58").Value = 1108
$ws.Range("I58").Value = 1111.8462
$ws.Range("K58").Value = 1111.8462
$ws.Range("M58").Value = -908.8462

$ws.Range("H132").Value = 6911.3335
$ws.Range("I132").Value = 8096
$ws.Range("K132").Value = 24288
$ws.Range("M132").Value = -21758

$ws.Range("H134").Value = 11905961
$ws.Range("I134").Value = 13890100
$ws.Range("J134").Value = 1128.5
$ws.Range("K134").Value = 41670300
$ws.Range("L134").Value = 3385.5
$ws.Range("M134").Value = -41667765
$ws.Range("N134").Value = -8455.5

$ws.Range("H136").Value = 1108
$ws.Range("I136").Value = 1111.8462
$ws.Range("K136").Value = 3335.5386
$ws.Range("M136").Value = -785.5385999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H104").Value = 4617.154
$ws.Range("I104").Value = 3341.6667
$ws.Range("K104").Value = 10025.0001
$ws.Range("M104").Value = -7404.000100000001

$ws.Range("H106").Value = 3165.3125
$ws.Range("J106").Value = 3165.3125
$ws.Range("L106").Value = 9495.9375
$ws.Range("N106").Value = -11387.9375

$ws.Range("H131").Value = 27028444
$ws.Range("J131").Value = 1751.2413
$ws.Range("L131").Value = 5253.7239
$ws.Range("N131").Value = -15333.7239

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 6200
$ws.Range("I29").Value = 6000
$ws.Range("K29").Value = 6000
$ws.Range("M29").Value = -5710

$ws.Range("H80").Value = 2741.5
$ws.Range("I80").Value = 1719.8
$ws.Range("J80").Value = 3471.2856
$ws.Range("K80").Value = 1719.8
$ws.Range("L80").Value = 3471.2856
$ws.Range("M80").Value = -721.8
$ws.Range("N80").Value = -5467.2856

$ws.Range("H83").Value = 2741.5
$ws.Range("I83").Value = 1719.8
$ws.Range("J83").Value = 3471.2856
$ws.Range("K83").Value = 8599
$ws.Range("L83").Value = 17356.428
$ws.Range("M83").Value = -3607
$ws.Range("N83").Value = -27340.428

$ws.Range("H102").Value = 5240.5
$ws.Range("I102").Value = 6520.6665
$ws.Range("K102").Value = 6520.6665
$ws.Range("M102").Value = -4898.6665

$ws.Range("H113").Value = 2688.9375
$ws.Range("I113").Value = 1366.6
$ws.Range("J113").Value = 3290
$ws.Range("K113").Value = 1366.6
$ws.Range("L113").Value = 3290
$ws.Range("M113").Value = 803.4000000000001
$ws.Range("N113").Value = -7630

$ws.Range("H132").Value = 1936.3334
$ws.Range("I132").Value = 1648.0555
$ws.Range("J132").Value = 3666
$ws.Range("K132").Value = 4944.166499999999
$ws.Range("L132").Value = 10998
$ws.Range("M132").Value = -2414.166499999999
$ws.Range("N132").Value = -16058

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3618.5
$ws.Range("I40").Value = 3289.6
$ws.Range("K40").Value = 3289.6
$ws.Range("M40").Value = -3153.6

$ws.Range("H96").Value = 30196.5
$ws.Range("J96").Value = 30196.5
$ws.Range("L96").Value = 30196.5
$ws.Range("N96").Value = -35688.5

$ws.Range("H132").Value = 79522
$ws.Range("I132").Value = 2198.111
$ws.Range("J132").Value = 253500.75
$ws.Range("K132").Value = 6594.333
$ws.Range("L132").Value = 760502.25
$ws.Range("M132").Value = -4064.333
$ws.Range("N132").Value = -765562.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4556.3687
$ws.Range("I132").Value = 3897.3333
$ws.Range("J132").Value = 5686.143
$ws.Range("K132").Value = 11691.9999
$ws.Range("L132").Value = 17058.429
$ws.Range("M132").Value = -9161.999899999999
$ws.Range("N132").Value = -22118.429
